$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sync the previous row's timestamp (refreshed source value with a
# negligible floating point shift vs. the prior snapshot)
$ws.Cells.Item(48, 1).Value2 = 44361.76756595949

$row = 49

$ws.Cells.Item($row, 1).Value2 = 44362.76909272074
$ws.Cells.Item($row, 2).Value = 77865
$ws.Cells.Item($row, 3).Value = 65438
$ws.Cells.Item($row, 4).Value = 3447
$ws.Cells.Item($row, 5).Value = 2102
$ws.Cells.Item($row, 6).Value = 1486
$ws.Cells.Item($row, 7).Value = 20578
$ws.Cells.Item($row, 8).Value = 1480
$ws.Cells.Item($row, 9).Value = 896
$ws.Cells.Item($row, 10).Value = 182

$ws.Range("A" + $row).NumberFormat = $ws.Range("A48").NumberFormat
